$wb = $excel.ActiveWorkbook

# --- "jobs" sheet: update Job Due Date values (column C) ---
$jobs = $wb.Worksheets.Item("jobs")
$jobs.Range("C2").Value = 43773.75
$jobs.Range("C3").Value = 43773.75
$jobs.Range("C4").Value = 43773.75
$jobs.Range("C5").Value = 43774.75
$jobs.Range("C6").Value = 43775.75
$jobs.Range("C7").Value = 43773.75
$jobs.Range("C8").Value = 43774.75
$jobs.Range("C9").Value = 43773.75
$jobs.Range("C10").Value = 43774.75
$jobs.Range("C11").Value = 43775.75

# --- "tasks" sheet: fix Task Runtime for row 36, and resize column F (best fit) ---
$tasks = $wb.Worksheets.Item("tasks")
$tasks.Range("D36").Value = 70
$tasks.Columns("F").ColumnWidth = 10.3
$tasks.Range("F3").Select()

# --- "machines" sheet: move selection ---
$machines = $wb.Worksheets.Item("machines")
$machines.Range("C6").Select()

# --- re-select "jobs" sheet last so it remains the active tab ---
$jobs.Range("G11").Select()
